# Insert a new weekly price record as row 386 in the daily-logic sheet,
# pushing the existing rows 386:407 down to 387:408.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 386:407 down by one to make room for the new record.
$ws.Rows("386:386").Insert()

# Copy the date cell's number formatting (style) from the row above so the
# new D386 renders the same way as every other date cell in the column.
$ws.Range("D385").Copy()
$ws.Range("D386").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row with the inserted record's data.
$ws.Range("A386").Value = 7
$ws.Range("B386").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C386").Value = "Ñuble"
$ws.Range("D386").Value = 45265
$ws.Range("E386").Value = 16
$ws.Range("F386").Value = 100112043
$ws.Range("G386").Value = "Pepino ensalada"
$ws.Range("H386").Value = "Sin especificar"
$ws.Range("I386").Value = "Primera"
$ws.Range("J386").Value = 120
$ws.Range("K386").Value = 18000
$ws.Range("L386").Value = 19000
$ws.Range("M386").Value = 18500
$ws.Range("N386").Value = "$/caja 80 unidades"
$ws.Range("O386").Value = "Región del Maule"
$ws.Range("P386").Value = 231
$ws.Range("Q386").Value = 80
$ws.Range("R386").Value = "Hortaliza"
